# Add "2022-Q3" data: new worksheet inserted right after "总计", plus a new
# summary row on "总计" itself.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (it becomes the
#    2nd tab, pushing "2022-Q2", "2022-Q1", ... one slot to the right).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @("014179", "中银证券远见价值混合A",       "1.56", "93.65", "7.87", "0.1228", 1),
    @("003980", "中银证券瑞益灵活配置混合A",     "0.66", "91.21", "6.82", "0.0450", 1),
    @("005571", "中银证券新能源灵活配置混合A",   "0.53", "90.32", "5.48", "0.0290", 6),
    @("164401", "前海开源中证健康产业指数",       "1.81", "94.19", "1.13", "0.0205", 9),
    @("005572", "中银证券新能源灵活配置混合C",   "0.25", "90.32", "5.48", "0.0137", 6),
    @("003981", "中银证券瑞益灵活配置混合C",     "0.19", "91.21", "6.82", "0.0130", 1),
    @("014180", "中银证券远见价值混合C",         "0.16", "93.65", "7.87", "0.0126", 1)
)

# Columns B (fund code) through G (holding value) are stored as text in the
# source workbook (even though several look numeric), so force text storage
# before writing, to avoid Excel auto-converting "014179" -> 14179.
$q3.Range("B2:G8").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $q3.Cells.Item($r, 1).Value = $i
    $q3.Cells.Item($r, 2).Value = $rows[$i][0]
    $q3.Cells.Item($r, 3).Value = $rows[$i][1]
    $q3.Cells.Item($r, 4).Value = $rows[$i][2]
    $q3.Cells.Item($r, 5).Value = $rows[$i][3]
    $q3.Cells.Item($r, 6).Value = $rows[$i][4]
    $q3.Cells.Item($r, 7).Value = $rows[$i][5]
    $q3.Cells.Item($r, 8).Value = $rows[$i][6]
}

# Match the "index column" / header styling (bold, centered, thin border)
# used throughout the rest of the workbook, by copying the format from the
# already-styled cells on the "总计" sheet.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$q3.Range("A2:A8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Update "总计": add a new "2022-Q3" row right after the header, and
#    shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q3", 7, 0.26),
    @("2022-Q2", 6, 0.26),
    @("2022-Q1", 1, 0.01),
    @("2021-Q4", 3, 0.11),
    @("2021-Q3", 7, 0.49),
    @("2021-Q2", 2, 0.01)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $summaryRows[$i][0]
    $summary.Cells.Item($r, 3).Value = $summaryRows[$i][1]
    $summary.Cells.Item($r, 4).Value = $summaryRows[$i][2]
}

# The new last row (row 7) needs the same index-column styling as the rest
# of column A.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore "总计" as the active sheet (matches the original workbook view).
$summary.Activate()
